# "Drop in results from RMI script"
#
# 1. Remove the "Texas Notes" sheet entirely.
# 2. Update the DR discount rate value (B2) from 5.87% to 3%.
# 3. Tidy up the selections left on the About / DR sheets.

$wb = $excel.ActiveWorkbook

# --- Remove the "Texas Notes" worksheet -----------------------------------
$notes = $wb.Worksheets.Item("Texas Notes")
$notes.Delete()

# --- Update the DR discount rate with the new RMI-script result -----------
$dr = $wb.Worksheets.Item("DR")
$dr.Range("B2").Value = 0.03

# --- Leave the cursor/selection where the author left it -------------------
$dr.Range("B1").Select()

$about = $wb.Worksheets.Item("About")
$about.Range("A16:A18").Select()

# Re-select the About sheet (it was the tab shown when the file was saved)
$about.Select()
